$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 2 (TS_01 / TC_01): Observed Result = same as Expected Result, Pass/Fail = Pass
$ws.Range("I2").Value = $ws.Range("H2").Value2
$ws.Range("J2").Value = "Pass"

# Row 6 (TS_02 / TC_01): Observed Result + Pass/Fail
$ws.Range("I6").Value = "All expenses are displayed in table"
$ws.Range("J6").Value = "Pass"

# Row 3 (TS_01 / TC_02): Expected Result updated with new validation message,
# Observed Result mirrors it, Pass/Fail = Pass
$ws.Range("H3").Value = 'Validation message shown, "Please enter amount."'
$ws.Range("I3").Value = $ws.Range("H3").Value2
$ws.Range("J3").Value = "Pass"

# Row 5: taller row to fit wrapped text
$ws.Rows.Item(5).RowHeight = 34.2

# Update the active selection to reflect where the author left off editing
$ws.Range("F4").Select()
